$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New row 36 should look like the other data rows (border + wrap formatting),
# so copy the formatting from the row directly above (row 35) first.
$ws.Range("A35:F35").Copy()
$ws.Range("A36:F36").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(36).RowHeight = 43.5

# Fill in the new test case data (TC035).
$ws.Range("A36").Value = "TC035"
$ws.Range("B36").Value = "Validation of the acknowledgment email sent to the email id"
$ws.Range("C36").Value = "1. Open the link in Browser`n2.Enter email field as 'saurabhsinghal001@gmail.com'`n3. Enter Submit Button"
$ws.Range("D36").Value = "High"
$ws.Range("E36").Value = "Email shold be delivered to the user"
$ws.Range("F36").Value = ""

# Move the visible selection down onto the freshly added row, like the
# author's cursor ended up after entering the new data.
$excel.ActiveWindow.ScrollRow = 26
$null = $ws.Range("A36:F36").Select()
